$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5, shifting existing rows 5-117 down to 6-118
$ws.Rows.Item(5).EntireRow.Insert()

# Populate the new row 5 with its data (same fixed fields as the rest of the
# dataset, plus the new observation's values)
$ws.Cells.Item(5, 1).Value = 11
$ws.Cells.Item(5, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(5, 3).Value = "Bíobío"
$ws.Cells.Item(5, 4).Value = 44817
$ws.Cells.Item(5, 5).Value = 8
$ws.Cells.Item(5, 6).Value = 100112001
$ws.Cells.Item(5, 7).Value = "Berenjena"
$ws.Cells.Item(5, 8).Value = "Sin especificar"
$ws.Cells.Item(5, 9).Value = "Primera"
$ws.Cells.Item(5, 10).Value = 180
$ws.Cells.Item(5, 11).Value = 12000
$ws.Cells.Item(5, 12).Value = 13000
$ws.Cells.Item(5, 13).Value = 12444
$ws.Cells.Item(5, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(5, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(5, 16).Value = 207
$ws.Cells.Item(5, 17).Value = 60
$ws.Cells.Item(5, 18).Value = "Hortaliza"

# D column uses a date/time number format on all data rows; make sure the
# newly-inserted row keeps it (Insert() already copies formatting from the
# row above, but set explicitly to be safe).
$ws.Cells.Item(5, 4).NumberFormat = $ws.Cells.Item(6, 4).NumberFormat
